$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, shifting existing rows 13:26 down to 14:27.
$ws.Rows("13:13").Insert()

# Populate the new row 13 with a new weekly record (same as the former row 13
# except for a new "Fecha" value), matching the inserted data row below it.
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C13").Value = "Ñuble"
$ws.Range("D13").Value = 44799
$ws.Range("D13").Style = $ws.Range("D14").Style
$ws.Range("D13").NumberFormat = $ws.Range("D14").NumberFormat
$ws.Range("E13").Value = 16
$ws.Range("F13").Value = 100112043
$ws.Range("G13").Value = "Pepino dulce"
$ws.Range("H13").Value = "Cultivar IV Región"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("N13").Value = "`$/bandeja 18 kilos"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 861
$ws.Range("Q13").Value = 18
$ws.Range("R13").Value = "Hortaliza"
